$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.451.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.097.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.089.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("E10").Value = "  +6.47%  "
$ws.Range("E11").Value = "  -3.38%  "
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.93%  "
$ws.Range("E15").Value = "  -1.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.611.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.281.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.090.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "459.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.721"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.09%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.110"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("E38").Value = "  -5.05%  "
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "431.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.877.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0366"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.268"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.75%  "
$ws.Range("E46").Value = "  -3.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.109"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.39%  "
